$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 3.560699333333333
$ws.Cells.Item(2, 8).Value = 10.682098
$ws.Cells.Item(2, 9).Value = 0.2516303646515017
$ws.Cells.Item(2, 10).Value = 0.2516303646515017
$ws.Cells.Item(2, 13).Value = 16.05260533333333
$ws.Cells.Item(2, 14).Value = 48.157816
$ws.Cells.Item(2, 15).Value = 0.1752915379534001
$ws.Cells.Item(2, 16).Value = 0.1752915379534001
$ws.Cells.Item(2, 17).Value = 57.15850110866311
$ws.Cells.Item(2, 18).Value = 514.426509977968
$ws.Cells.Item(2, 19).Value = 0.04410867361553662
$ws.Cells.Item(2, 20).Value = 0.04410867361553662

# Row 3
$ws.Cells.Item(3, 7).Value = 3.560699333333333
$ws.Cells.Item(3, 8).Value = 10.682098
$ws.Cells.Item(3, 9).Value = 0.2516303646515017
$ws.Cells.Item(3, 10).Value = 0.2516303646515017
$ws.Cells.Item(3, 15).Value = 0.07888758308485012
$ws.Cells.Item(3, 16).Value = 0.07888758308485012
$ws.Cells.Item(3, 17).Value = 25.72340945752822
$ws.Cells.Item(3, 18).Value = 231.510685117754
$ws.Cells.Item(3, 19).Value = 0.01985051129811648
$ws.Cells.Item(3, 20).Value = 0.01985051129811648

# Row 4
$ws.Cells.Item(4, 7).Value = 3.560699333333333
$ws.Cells.Item(4, 8).Value = 10.682098
$ws.Cells.Item(4, 9).Value = 0.2516303646515017
$ws.Cells.Item(4, 10).Value = 0.2516303646515017
$ws.Cells.Item(4, 13).Value = 2.098187333333334
$ws.Cells.Item(4, 14).Value = 6.294562000000001
$ws.Cells.Item(4, 15).Value = 0.02291182502385553
$ws.Cells.Item(4, 16).Value = 0.02291182502385553
$ws.Cells.Item(4, 17).Value = 7.471014239008445
$ws.Cells.Item(4, 18).Value = 67.239128151076
$ws.Cells.Item(4, 19).Value = 0.005765310885584169
$ws.Cells.Item(4, 20).Value = 0.005765310885584169

# Row 5
$ws.Cells.Item(5, 7).Value = 3.560699333333333
$ws.Cells.Item(5, 8).Value = 10.682098
$ws.Cells.Item(5, 9).Value = 0.2516303646515017
$ws.Cells.Item(5, 10).Value = 0.2516303646515017
$ws.Cells.Item(5, 13).Value = 66.20156266666667
$ws.Cells.Item(5, 14).Value = 198.604688
$ws.Cells.Item(5, 15).Value = 0.7229090539378943
$ws.Cells.Item(5, 16).Value = 0.7229090539378942
$ws.Cells.Item(5, 17).Value = 235.7238600528249
$ws.Cells.Item(5, 18).Value = 2121.514740475424
$ws.Cells.Item(5, 19).Value = 0.1819058688522645
$ws.Cells.Item(5, 20).Value = 0.1819058688522645

# Row 6
$ws.Cells.Item(6, 9).Value = 0.2153092375010323
$ws.Cells.Item(6, 10).Value = 0.2153092375010323
$ws.Cells.Item(6, 13).Value = 16.05260533333333
$ws.Cells.Item(6, 14).Value = 48.157816
$ws.Cells.Item(6, 15).Value = 0.1752915379534001
$ws.Cells.Item(6, 16).Value = 0.1752915379534001
$ws.Cells.Item(6, 17).Value = 48.90806126459555
$ws.Cells.Item(6, 18).Value = 440.1725513813599
$ws.Cells.Item(6, 19).Value = 0.03774188737712984
$ws.Cells.Item(6, 20).Value = 0.03774188737712984

# Row 7
$ws.Cells.Item(7, 9).Value = 0.2153092375010323
$ws.Cells.Item(7, 10).Value = 0.2153092375010323
$ws.Cells.Item(7, 15).Value = 0.07888758308485012
$ws.Cells.Item(7, 16).Value = 0.07888758308485012
$ws.Cells.Item(7, 19).Value = 0.01698522536229841
$ws.Cells.Item(7, 20).Value = 0.01698522536229842

# Row 8
$ws.Cells.Item(8, 9).Value = 0.2153092375010323
$ws.Cells.Item(8, 10).Value = 0.2153092375010323
$ws.Cells.Item(8, 13).Value = 2.098187333333334
$ws.Cells.Item(8, 14).Value = 6.294562000000001
$ws.Cells.Item(8, 15).Value = 0.02291182502385553
$ws.Cells.Item(8, 16).Value = 0.02291182502385553
$ws.Cells.Item(8, 17).Value = 6.392624282002223
$ws.Cells.Item(8, 18).Value = 57.53361853802001
$ws.Cells.Item(8, 19).Value = 0.004933127575643405
$ws.Cells.Item(8, 20).Value = 0.004933127575643406

# Row 9
$ws.Cells.Item(9, 9).Value = 0.2153092375010323
$ws.Cells.Item(9, 10).Value = 0.2153092375010323
$ws.Cells.Item(9, 13).Value = 66.20156266666667
$ws.Cells.Item(9, 14).Value = 198.604688
$ws.Cells.Item(9, 15).Value = 0.7229090539378943
$ws.Cells.Item(9, 16).Value = 0.7229090539378942
$ws.Cells.Item(9, 17).Value = 201.6987283671644
$ws.Cells.Item(9, 18).Value = 1815.28855530448
$ws.Cells.Item(9, 19).Value = 0.1556489971859606
$ws.Cells.Item(9, 20).Value = 0.1556489971859606

# Row 10
$ws.Cells.Item(10, 7).Value = 4.835201333333333
$ws.Cells.Item(10, 8).Value = 14.505604
$ws.Cells.Item(10, 9).Value = 0.3416978971743455
$ws.Cells.Item(10, 10).Value = 0.3416978971743456
$ws.Cells.Item(10, 13).Value = 16.05260533333333
$ws.Cells.Item(10, 14).Value = 48.157816
$ws.Cells.Item(10, 15).Value = 0.1752915379534001
$ws.Cells.Item(10, 16).Value = 0.1752915379534001
$ws.Cells.Item(10, 17).Value = 77.6175787112071
$ws.Cells.Item(10, 18).Value = 698.5582084008639
$ws.Cells.Item(10, 19).Value = 0.05989674991113379
$ws.Cells.Item(10, 20).Value = 0.05989674991113379

# Row 11
$ws.Cells.Item(11, 7).Value = 4.835201333333333
$ws.Cells.Item(11, 8).Value = 14.505604
$ws.Cells.Item(11, 9).Value = 0.3416978971743455
$ws.Cells.Item(11, 10).Value = 0.3416978971743456
$ws.Cells.Item(11, 15).Value = 0.07888758308485012
$ws.Cells.Item(11, 16).Value = 0.07888758308485012
$ws.Cells.Item(11, 17).Value = 34.93074030221022
$ws.Cells.Item(11, 18).Value = 314.376662719892
$ws.Cells.Item(11, 19).Value = 0.02695572125325976
$ws.Cells.Item(11, 20).Value = 0.02695572125325976

# Row 12
$ws.Cells.Item(12, 7).Value = 4.835201333333333
$ws.Cells.Item(12, 8).Value = 14.505604
$ws.Cells.Item(12, 9).Value = 0.3416978971743455
$ws.Cells.Item(12, 10).Value = 0.3416978971743456
$ws.Cells.Item(12, 13).Value = 2.098187333333334
$ws.Cells.Item(12, 14).Value = 6.294562000000001
$ws.Cells.Item(12, 15).Value = 0.02291182502385553
$ws.Cells.Item(12, 16).Value = 0.02291182502385553
$ws.Cells.Item(12, 17).Value = 10.14515819171645
$ws.Cells.Item(12, 18).Value = 91.30642372544801
$ws.Cells.Item(12, 19).Value = 0.007828922431077983
$ws.Cells.Item(12, 20).Value = 0.007828922431077984

# Row 13
$ws.Cells.Item(13, 7).Value = 4.835201333333333
$ws.Cells.Item(13, 8).Value = 14.505604
$ws.Cells.Item(13, 9).Value = 0.3416978971743455
$ws.Cells.Item(13, 10).Value = 0.3416978971743456
$ws.Cells.Item(13, 13).Value = 66.20156266666667
$ws.Cells.Item(13, 14).Value = 198.604688
$ws.Cells.Item(13, 15).Value = 0.7229090539378943
$ws.Cells.Item(13, 16).Value = 0.7229090539378942
$ws.Cells.Item(13, 17).Value = 320.0978840746169
$ws.Cells.Item(13, 18).Value = 2880.880956671552
$ws.Cells.Item(13, 19).Value = 0.247016503578874
$ws.Cells.Item(13, 20).Value = 0.247016503578874

# Row 14
$ws.Cells.Item(14, 7).Value = 2.707878
$ws.Cells.Item(14, 8).Value = 8.123634000000001
$ws.Cells.Item(14, 9).Value = 0.1913625006731204
$ws.Cells.Item(14, 10).Value = 0.1913625006731204
$ws.Cells.Item(14, 13).Value = 16.05260533333333
$ws.Cells.Item(14, 14).Value = 48.157816
$ws.Cells.Item(14, 15).Value = 0.1752915379534001
$ws.Cells.Item(14, 16).Value = 0.1752915379534001
$ws.Cells.Item(14, 17).Value = 43.46849682481601
$ws.Cells.Item(14, 18).Value = 391.216471423344
$ws.Cells.Item(14, 19).Value = 0.03354422704959983
$ws.Cells.Item(14, 20).Value = 0.03354422704959983

# Row 15
$ws.Cells.Item(15, 7).Value = 2.707878
$ws.Cells.Item(15, 8).Value = 8.123634000000001
$ws.Cells.Item(15, 9).Value = 0.1913625006731204
$ws.Cells.Item(15, 10).Value = 0.1913625006731204
$ws.Cells.Item(15, 15).Value = 0.07888758308485012
$ws.Cells.Item(15, 16).Value = 0.07888758308485012
$ws.Cells.Item(15, 17).Value = 19.562408401898
$ws.Cells.Item(15, 18).Value = 176.061675617082
$ws.Cells.Item(15, 19).Value = 0.01509612517117547
$ws.Cells.Item(15, 20).Value = 0.01509612517117547

# Row 16
$ws.Cells.Item(16, 7).Value = 2.707878
$ws.Cells.Item(16, 8).Value = 8.123634000000001
$ws.Cells.Item(16, 9).Value = 0.1913625006731204
$ws.Cells.Item(16, 10).Value = 0.1913625006731204
$ws.Cells.Item(16, 13).Value = 2.098187333333334
$ws.Cells.Item(16, 14).Value = 6.294562000000001
$ws.Cells.Item(16, 15).Value = 0.02291182502385553
$ws.Cells.Item(16, 16).Value = 0.02291182502385553
$ws.Cells.Item(16, 17).Value = 5.681635319812002
$ws.Cells.Item(16, 18).Value = 51.13471787830801
$ws.Cells.Item(16, 19).Value = 0.00438446413154997
$ws.Cells.Item(16, 20).Value = 0.00438446413154997

# Row 17
$ws.Cells.Item(17, 7).Value = 2.707878
$ws.Cells.Item(17, 8).Value = 8.123634000000001
$ws.Cells.Item(17, 9).Value = 0.1913625006731204
$ws.Cells.Item(17, 10).Value = 0.1913625006731204
$ws.Cells.Item(17, 13).Value = 66.20156266666667
$ws.Cells.Item(17, 14).Value = 198.604688
$ws.Cells.Item(17, 15).Value = 0.7229090539378943
$ws.Cells.Item(17, 16).Value = 0.7229090539378942
$ws.Cells.Item(17, 17).Value = 179.265755110688
$ws.Cells.Item(17, 18).Value = 1613.391795996192
$ws.Cells.Item(17, 19).Value = 0.1383376843207951
$ws.Cells.Item(17, 20).Value = 0.1383376843207951
